$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 210, pushing existing rows 210-333 down to 212-335.
$ws.Rows.Item(210).Insert()
$ws.Rows.Item(210).Insert()

# Fill in the new row 210 (Primera) with the new weekly data point.
$ws.Range("A210").Value = 11
$ws.Range("B210").Value = "Vega Monumental Concepción"
$ws.Range("C210").Value = "Bíobío"
$ws.Range("D210").Value = 44806
$ws.Range("E210").Value = 8
$ws.Range("F210").Value = 100112017
$ws.Range("G210").Value = "Apio"
$ws.Range("H210").Value = "Americana (o)"
$ws.Range("I210").Value = "Primera"
$ws.Range("J210").Value = 220
$ws.Range("K210").Value = 8000
$ws.Range("L210").Value = 8500
$ws.Range("M210").Value = 8227
$ws.Range("N210").Value = "`$/docena de matas"
$ws.Range("O210").Value = "Región de Coquimbo"
$ws.Range("P210").Value = 1371
$ws.Range("Q210").Value = 6
$ws.Range("R210").Value = "Hortaliza"

# Fill in the new row 211 (Segunda) with the new weekly data point.
$ws.Range("A211").Value = 11
$ws.Range("B211").Value = "Vega Monumental Concepción"
$ws.Range("C211").Value = "Bíobío"
$ws.Range("D211").Value = 44806
$ws.Range("E211").Value = 8
$ws.Range("F211").Value = 100112017
$ws.Range("G211").Value = "Apio"
$ws.Range("H211").Value = "Americana (o)"
$ws.Range("I211").Value = "Segunda"
$ws.Range("J211").Value = 220
$ws.Range("K211").Value = 6500
$ws.Range("L211").Value = 7000
$ws.Range("M211").Value = 6727
$ws.Range("N211").Value = "`$/docena de matas"
$ws.Range("O211").Value = "Región de Coquimbo"
$ws.Range("P211").Value = 1121
$ws.Range("Q211").Value = 6
$ws.Range("R211").Value = "Hortaliza"

# Match the date formatting style used by the rest of column D.
$ws.Range("D210").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D211").NumberFormat = "YYYY-MM-DD HH:MM:SS"
